$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Mutombo" / "Bongo" entries (currently A18 and A19) ---
$ws.Rows("18:19").Delete() | Out-Null

# --- Close up the blank-row gaps that existed between groups of slurs ---
# After removing rows 18:19 above, the rows below shift up by two, so the
# blank rows that used to sit at A10, A16 and A21 are now at A10, A16, A19.
# Delete them bottom-up so earlier row numbers stay valid.
$ws.Rows("19:19").Delete() | Out-Null
$ws.Rows("16:16").Delete() | Out-Null
$ws.Rows("10:10").Delete() | Out-Null

# --- Update the view: scroll so row 4 is at the top, and select B12 ---
$ws.Range("B12").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
